$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = -0.734317941686569
$ws.Range("F2").Value = 0.0908559981901698
$ws.Range("G2").Value = -8.08221753449426
$ws.Range("H2").Value = 0.000000000000000635995242796327

# Row 3
$ws.Range("E3").Value = 0.430321032504187
$ws.Range("F3").Value = 0.0528291680891282
$ws.Range("G3").Value = 8.14551975867178
$ws.Range("H3").Value = 0.000000000000000377656166135215

# Row 4
$ws.Range("E4").Value = 0.477007183896796
$ws.Range("F4").Value = 0.0522120478110245
$ws.Range("G4").Value = 9.13596006851271
$ws.Range("H4").Value = 0.0000000000000000000648289627606608

# Row 5
$ws.Range("E5").Value = 0.101955562126693
$ws.Range("F5").Value = 0.0500936631602483
$ws.Range("G5").Value = 2.03529859256928
$ws.Range("H5").Value = 0.041820837711279

# Row 6
$ws.Range("E6").Value = -0.0576262458677776
$ws.Range("F6").Value = 0.00919657263136726
$ws.Range("G6").Value = -6.26605673413905
$ws.Range("H6").Value = 0.000000000370305464104952

# Row 7
$ws.Range("E7").Value = 0.355007925948079
